# Daily attendance processing - 2025-12-21 04:24:32
# Re-order the "Recorded By" (column G) list so that "System" (exact case)
# is moved to the front of the comma-separated list, when present but not
# already first. Other entries (including a duplicate lowercase "system",
# and email addresses) keep their existing relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val.ToString().Contains("System")) {
        $parts = $val.ToString().Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        if ($trimmed[0] -ne "System") {
            $rest = @()
            $found = $false
            foreach ($p in $trimmed) {
                if (-not $found -and $p -eq "System") {
                    $found = $true
                } else {
                    $rest += $p
                }
            }

            if ($found) {
                $newParts = @("System") + $rest
                $newVal = [string]::Join(", ", $newParts)
                $cell.Value2 = $newVal
            }
        }
    }
}
